# Scheduled-runner update: refresh Leve profit-calc columns (H-N) across all
# crafting-class sheets with the latest currentAveragePrice market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1876.8334
$ws.Range("I40").Value = 1852.2
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1852.2
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1677.2
$ws.Range("N40").Value = -2350

$ws.Range("H51").Value = 5373.5
$ws.Range("J51").Value = 5373.5
$ws.Range("L51").Value = 5373.5
$ws.Range("N51").Value = -6341.5

$ws.Range("H62").Value = 1563.8049
$ws.Range("I62").Value = 1633.2142
$ws.Range("J62").Value = 1414.3077
$ws.Range("K62").Value = 1633.2142
$ws.Range("L62").Value = 1414.3077
$ws.Range("M62").Value = -1009.2142
$ws.Range("N62").Value = -2662.3077

$ws.Range("H65").Value = 1563.8049
$ws.Range("I65").Value = 1633.2142
$ws.Range("J65").Value = 1414.3077
$ws.Range("K65").Value = 8166.071
$ws.Range("L65").Value = 7071.538500000001
$ws.Range("M65").Value = -5046.071
$ws.Range("N65").Value = -13311.5385

$ws.Range("H100").Value = 16668780
$ws.Range("I100").Value = 23810900
$ws.Range("K100").Value = 23810900
$ws.Range("M100").Value = -23810359

$ws.Range("H116").Value = 12630.5
$ws.Range("I116").Value = 21461
$ws.Range("J116").Value = 3800
$ws.Range("K116").Value = 21461
$ws.Range("L116").Value = 3800
$ws.Range("M116").Value = -18019
$ws.Range("N116").Value = -10684

$ws.Range("H132").Value = 47082.59
$ws.Range("I132").Value = 53579.844
$ws.Range("J132").Value = 5933.3335
$ws.Range("K132").Value = 160739.532
$ws.Range("L132").Value = 17800.0005
$ws.Range("M132").Value = -158209.532
$ws.Range("N132").Value = -22860.0005

$ws.Range("H133").Value = 59800
$ws.Range("J133").Value = 59800
$ws.Range("L133").Value = 59800
$ws.Range("N133").Value = -69920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 528392
$ws.Range("I32").Value = 5033.5776
$ws.Range("J32").Value = 2883505
$ws.Range("K32").Value = 5033.5776
$ws.Range("L32").Value = 2883505
$ws.Range("M32").Value = -4746.5776
$ws.Range("N32").Value = -2884079

$ws.Range("H63").Value = 4448
$ws.Range("I63").Value = 3151.5557
$ws.Range("J63").Value = 5744.4443
$ws.Range("K63").Value = 3151.5557
$ws.Range("L63").Value = 5744.4443
$ws.Range("M63").Value = -2465.5557
$ws.Range("N63").Value = -7116.4443

$ws.Range("H66").Value = 4448
$ws.Range("I66").Value = 3151.5557
$ws.Range("J66").Value = 5744.4443
$ws.Range("K66").Value = 15757.7785
$ws.Range("L66").Value = 28722.2215
$ws.Range("M66").Value = -12325.7785
$ws.Range("N66").Value = -35586.2215

$ws.Range("H122").Value = 62134.234
$ws.Range("I122").Value = 3302.8572
$ws.Range("K122").Value = 9908.571599999999
$ws.Range("M122").Value = -7458.571599999999

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 438.7
$ws.Range("I80").Value = 697.44446
$ws.Range("J80").Value = 227
$ws.Range("K80").Value = 697.44446
$ws.Range("L80").Value = 227
$ws.Range("M80").Value = 300.55554
$ws.Range("N80").Value = -2223

$ws.Range("H83").Value = 438.7
$ws.Range("I83").Value = 697.44446
$ws.Range("J83").Value = 227
$ws.Range("K83").Value = 3487.2223
$ws.Range("L83").Value = 1135
$ws.Range("M83").Value = 1504.7777
$ws.Range("N83").Value = -11119

$ws.Range("H94").Value = 1156.125
$ws.Range("I94").Value = 1423.1666
$ws.Range("K94").Value = 1423.1666
$ws.Range("M94").Value = -972.1666

$ws.Range("H99").Value = 3370.1428
$ws.Range("I99").Value = 2990
$ws.Range("J99").Value = 3522.2
$ws.Range("K99").Value = 2990
$ws.Range("L99").Value = 3522.2
$ws.Range("M99").Value = -1492
$ws.Range("N99").Value = -6518.2

$ws.Range("H132").Value = 19000
$ws.Range("J132").Value = 19000
$ws.Range("L132").Value = 19000
$ws.Range("N132").Value = -29120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6136.615
$ws.Range("I31").Value = 6778
$ws.Range("K31").Value = 6778
$ws.Range("M31").Value = -6483

$ws.Range("H34").Value = 6136.615
$ws.Range("I34").Value = 6778
$ws.Range("K34").Value = 6778
$ws.Range("M34").Value = -6576

$ws.Range("H62").Value = 1533.3334
$ws.Range("I62").Value = 1300
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1300
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -676
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 1533.3334
$ws.Range("I65").Value = 1300
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 6500
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -3380
$ws.Range("N65").Value = -16240

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H104").Value = 49990
$ws.Range("J104").Value = 49990
$ws.Range("L104").Value = 49990
$ws.Range("N104").Value = -55232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 739.88464
$ws.Range("I5").Value = 653.8261
$ws.Range("K5").Value = 1961.4783
$ws.Range("M5").Value = -1849.4783

$ws.Range("H118").Value = 2786.111
$ws.Range("I118").Value = 1950
$ws.Range("J118").Value = 2835.2942
$ws.Range("K118").Value = 5850
$ws.Range("L118").Value = 8505.882599999999
$ws.Range("M118").Value = -4607
$ws.Range("N118").Value = -10991.8826

$ws.Range("H135").Value = 739.88464
$ws.Range("I135").Value = 653.8261
$ws.Range("K135").Value = 5884.4349
$ws.Range("M135").Value = -3349.4349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 188334.83
$ws.Range("I80").Value = 556505
$ws.Range("J80").Value = 4249.75
$ws.Range("K80").Value = 556505
$ws.Range("L80").Value = 4249.75
$ws.Range("M80").Value = -555507
$ws.Range("N80").Value = -6245.75

$ws.Range("H83").Value = 188334.83
$ws.Range("I83").Value = 556505
$ws.Range("J83").Value = 4249.75
$ws.Range("K83").Value = 2782525
$ws.Range("L83").Value = 21248.75
$ws.Range("M83").Value = -2777533
$ws.Range("N83").Value = -31232.75

$ws.Range("H97").Value = 1376.6666
$ws.Range("I97").Value = 708.6
$ws.Range("J97").Value = 2490.111
$ws.Range("K97").Value = 708.6
$ws.Range("L97").Value = 2490.111
$ws.Range("M97").Value = -212.6
$ws.Range("N97").Value = -3482.111

$ws.Range("H130").Value = 54955
$ws.Range("J130").Value = 54955
$ws.Range("L130").Value = 54955
$ws.Range("N130").Value = -64995

$ws.Range("H133").Value = 42813.332
$ws.Range("J133").Value = 42813.332
$ws.Range("L133").Value = 42813.332
$ws.Range("N133").Value = -52933.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1256.5
$ws.Range("I32").Value = 707.8
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 707.8
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -390.8
$ws.Range("N32").Value = -4634

$ws.Range("H40").Value = 1927.5938
$ws.Range("I40").Value = 1776.3529
$ws.Range("J40").Value = 2099
$ws.Range("K40").Value = 1776.3529
$ws.Range("L40").Value = 2099
$ws.Range("M40").Value = -1640.3529
$ws.Range("N40").Value = -2371

$ws.Range("H46").Value = 1344.0555
$ws.Range("I46").Value = 1332.75
$ws.Range("J46").Value = 1366.6666
$ws.Range("K46").Value = 1332.75
$ws.Range("L46").Value = 1366.6666
$ws.Range("M46").Value = -1144.75
$ws.Range("N46").Value = -1742.6666

$ws.Range("H68").Value = 2057.7354
$ws.Range("I68").Value = 1961.8182
$ws.Range("J68").Value = 2233.5833
$ws.Range("K68").Value = 1961.8182
$ws.Range("L68").Value = 2233.5833
$ws.Range("M68").Value = -1212.8182
$ws.Range("N68").Value = -3731.5833

$ws.Range("H71").Value = 2057.7354
$ws.Range("I71").Value = 1961.8182
$ws.Range("J71").Value = 2233.5833
$ws.Range("K71").Value = 9809.091
$ws.Range("L71").Value = 11167.9165
$ws.Range("M71").Value = -6065.091
$ws.Range("N71").Value = -18655.9165

$ws.Range("H82").Value = 1992
$ws.Range("I82").Value = 1777.6666
$ws.Range("J82").Value = 2359.4285
$ws.Range("K82").Value = 1777.6666
$ws.Range("L82").Value = 2359.4285
$ws.Range("M82").Value = -1416.6666
$ws.Range("N82").Value = -3081.4285

$ws.Range("H85").Value = 1992
$ws.Range("I85").Value = 1777.6666
$ws.Range("J85").Value = 2359.4285
$ws.Range("K85").Value = 1777.6666
$ws.Range("L85").Value = 2359.4285
$ws.Range("M85").Value = -529.6666
$ws.Range("N85").Value = -4855.4285

$ws.Range("H133").Value = 84600
$ws.Range("J133").Value = 84600
$ws.Range("L133").Value = 84600
$ws.Range("N133").Value = -89660

$ws.Range("H139").Value = 70533.336
$ws.Range("J139").Value = 70533.336
$ws.Range("L139").Value = 70533.336
$ws.Range("N139").Value = -80813.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 34715
$ws.Range("J133").Value = 34715
$ws.Range("L133").Value = 34715
$ws.Range("N133").Value = -44835
